$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B4").Value = 7.645999999999999
$ws.Range("B6").Value = 6.944999999999999
$ws.Range("B7").Value = 5.05
$ws.Range("C7").Value = -13.601
$ws.Range("C12").Value = -11.131
$ws.Range("D13").Value = -8.007999999999999
$ws.Range("D14").Value = -7.774000000000001
$ws.Range("C15").Value = -13.609
$ws.Range("B16").Value = 5.524999999999999
$ws.Range("D16").Value = -8.614000000000001
$ws.Range("D19").Value = -8.184999999999999
$ws.Range("B20").Value = 8.888999999999999
$ws.Range("C20").Value = -12.148
$ws.Range("C21").Value = -12.182
$ws.Range("C22").Value = -12.925
$ws.Range("D22").Value = -7.818000000000001
$ws.Range("C23").Value = -12.451
$ws.Range("B28").Value = 6.641
$ws.Range("B29").Value = 5.255
$ws.Range("C29").Value = -11.358
$ws.Range("B32").Value = 6.431
$ws.Range("C34").Value = -12.44
$ws.Range("D36").Value = -7.987
$ws.Range("B40").Value = 9.263999999999999
$ws.Range("C42").Value = -11.999
$ws.Range("C43").Value = -13.867
$ws.Range("C44").Value = -13.636
$ws.Range("C45").Value = -13.376
$ws.Range("B46").Value = 5.545
$ws.Range("C46").Value = -14.141
$ws.Range("D46").Value = -8.370000000000001
$ws.Range("C50").Value = -13.625
$ws.Range("D50").Value = -8.518000000000001
$ws.Range("B51").Value = 5.059
$ws.Range("C51").Value = -12.057
$ws.Range("B52").Value = 5.442
$ws.Range("B57").Value = 5.763
$ws.Range("B59").Value = 5.329000000000001
$ws.Range("B62").Value = 5.902
$ws.Range("B66").Value = 4.961
$ws.Range("C66").Value = -10.897
$ws.Range("C67").Value = -11.3
$ws.Range("B73").Value = 7.316
$ws.Range("B74").Value = 9.132999999999999
$ws.Range("C79").Value = -12.204
$ws.Range("C84").Value = -13.68
$ws.Range("B92").Value = 6.272
$ws.Range("C92").Value = -10.898
$ws.Range("D95").Value = -7.753
$ws.Range("C97").Value = -12.146
$ws.Range("D97").Value = -8.599
$ws.Range("B100").Value = 6.486
